# Generate Report for Handback
# Updates the localization-status workbook: marks all "Ready for handoff"
# rows as handed back, refreshes the handback timestamps, and clears the
# stale "version not latest" error detail now that the handback is in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.14437166849777
$wsOverview.Columns.Item(6).ColumnWidth = 29.14437166849777

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-10-17 17:14:30"
$wsZhCn.Range("K3").Value = "2016-10-17 17:14:30"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsZhCn.Columns.Item(16).ColumnWidth = 12.913719813028965

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-10-17 17:15:14"
$wsDeDe.Range("K3").Value = "2016-10-17 17:15:14"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsDeDe.Columns.Item(16).ColumnWidth = 12.913719813028965
